# "Ran 1st analysis in Python."
#
# The raw Titanic export living on "Sheet1" is the table the analysis
# actually works from, so rename it to something more descriptive ("data")
# and switch focus to it (it was "Sheet2" - the PivotTable summary sheet -
# that had the selection/focus before).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "data"

# Activating the sheet makes it the workbook's active tab (bookViews/
# activeTab) and moves tabSelected from the old sheet to this one.
$ws.Activate()
